$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.796007752418518
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 2.527757883071899
$ws.Range("D1").Value = 1.100344181060791
$ws.Range("E1").Value = 0.734033465385437
